$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(2, 8).Value = 223.46666
$ws.Cells.Item(2, 9).Value = 180
$ws.Cells.Item(2, 10).Value = 506
$ws.Cells.Item(2, 11).Value = 180
$ws.Cells.Item(2, 12).Value = 506
$ws.Cells.Item(2, 13).Value = -67
$ws.Cells.Item(2, 14).Value = -732
$ws.Cells.Item(6, 8).Value = 1109.8
$ws.Cells.Item(6, 9).Value = 1362.25
$ws.Cells.Item(6, 11).Value = 4086.75
$ws.Cells.Item(6, 13).Value = -3974.75
$ws.Cells.Item(8, 8).Value = 1234.4445
$ws.Cells.Item(8, 9).Value = 1234.4445
$ws.Cells.Item(8, 11).Value = 3703.3335
$ws.Cells.Item(8, 13).Value = -3564.3335
$ws.Cells.Item(21, 8).Value = 0
$ws.Cells.Item(21, 9).Value = 0
$ws.Cells.Item(21, 11).Value = 0
$ws.Cells.Item(21, 13).ClearContents()
$ws.Cells.Item(23, 8).Value = 0
$ws.Cells.Item(23, 9).Value = 0
$ws.Cells.Item(23, 11).Value = 0
$ws.Cells.Item(23, 13).ClearContents()
$ws.Cells.Item(34, 8).Value = 10616
$ws.Cells.Item(34, 9).Value = 10616
$ws.Cells.Item(34, 11).Value = 10616
$ws.Cells.Item(34, 13).Value = -10413
$ws.Cells.Item(36, 8).Value = 10616
$ws.Cells.Item(36, 9).Value = 10616
$ws.Cells.Item(36, 11).Value = 10616
$ws.Cells.Item(36, 13).Value = -9901
$ws.Cells.Item(38, 8).Value = 2868
$ws.Cells.Item(38, 9).Value = 1705
$ws.Cells.Item(38, 11).Value = 5115
$ws.Cells.Item(38, 13).Value = -4743
$ws.Cells.Item(39, 8).Value = 157.55556
$ws.Cells.Item(39, 9).Value = 85.85714
$ws.Cells.Item(39, 11).Value = 257.57142
$ws.Cells.Item(39, 13).Value = 38.42858000000001
$ws.Cells.Item(40, 8).Value = 3220.8667
$ws.Cells.Item(40, 9).Value = 2404.4062
$ws.Cells.Item(40, 11).Value = 2404.4062
$ws.Cells.Item(40, 13).Value = -2229.4062
$ws.Cells.Item(62, 8).Value = 2849.6428
$ws.Cells.Item(62, 9).Value = 2849.6428
$ws.Cells.Item(62, 11).Value = 2849.6428
$ws.Cells.Item(62, 13).Value = -2225.6428
$ws.Cells.Item(65, 8).Value = 2849.6428
$ws.Cells.Item(65, 9).Value = 2849.6428
$ws.Cells.Item(65, 11).Value = 14248.214
$ws.Cells.Item(65, 13).Value = -11128.214
$ws.Cells.Item(87, 8).Value = 87499.5
$ws.Cells.Item(87, 9).Value = 49999.668
$ws.Cells.Item(87, 10).Value = 199999
$ws.Cells.Item(87, 11).Value = 49999.668
$ws.Cells.Item(87, 12).Value = 199999
$ws.Cells.Item(87, 13).Value = -48751.668
$ws.Cells.Item(87, 14).Value = -202495
$ws.Cells.Item(90, 8).Value = 87499.5
$ws.Cells.Item(90, 9).Value = 49999.668
$ws.Cells.Item(90, 10).Value = 199999
$ws.Cells.Item(90, 11).Value = 149999.004
$ws.Cells.Item(90, 12).Value = 599997
$ws.Cells.Item(90, 13).Value = -143759.004
$ws.Cells.Item(90, 14).Value = -612477
$ws.Cells.Item(92, 8).Value = 1191.8235
$ws.Cells.Item(92, 9).Value = 1217.4
$ws.Cells.Item(92, 11).Value = 1217.4
$ws.Cells.Item(92, 13).Value = 30.59999999999991
$ws.Cells.Item(99, 8).Value = 1079
$ws.Cells.Item(99, 9).Value = 191.44444
$ws.Cells.Item(99, 10).Value = 5073
$ws.Cells.Item(99, 11).Value = 574.33332
$ws.Cells.Item(99, 12).Value = 15219
$ws.Cells.Item(99, 13).Value = 923.66668
$ws.Cells.Item(99, 14).Value = -18215
$ws.Cells.Item(107, 8).Value = 1125.5
$ws.Cells.Item(107, 10).Value = 1300.75
$ws.Cells.Item(107, 12).Value = 1300.75
$ws.Cells.Item(107, 14).Value = -5140.75
$ws.Cells.Item(129, 8).Value = 2039.6364
$ws.Cells.Item(129, 9).Value = 990.1667
$ws.Cells.Item(129, 10).Value = 3299
$ws.Cells.Item(129, 11).Value = 2970.5001
$ws.Cells.Item(129, 12).Value = 9897
$ws.Cells.Item(129, 13).Value = 2029.4999
$ws.Cells.Item(129, 14).Value = -19897
$ws.Cells.Item(132, 8).Value = 5275.3447
$ws.Cells.Item(132, 9).Value = 4895.926
$ws.Cells.Item(132, 10).Value = 10397.5
$ws.Cells.Item(132, 11).Value = 14687.778
$ws.Cells.Item(132, 12).Value = 31192.5
$ws.Cells.Item(132, 13).Value = -12157.778
$ws.Cells.Item(132, 14).Value = -36252.5
$ws.Cells.Item(133, 8).Value = 107993.5
$ws.Cells.Item(133, 10).Value = 107993.5
$ws.Cells.Item(133, 12).Value = 107993.5
$ws.Cells.Item(133, 14).Value = -118113.5
$ws.Cells.Item(135, 8).Value = 823.36365
$ws.Cells.Item(135, 9).Value = 817.7143
$ws.Cells.Item(135, 10).Value = 833.25
$ws.Cells.Item(135, 11).Value = 7359.428699999999
$ws.Cells.Item(135, 12).Value = 7499.25
$ws.Cells.Item(135, 13).Value = -4824.428699999999
$ws.Cells.Item(135, 14).Value = -12569.25
$ws.Cells.Item(136, 8).Value = 109990.336
$ws.Cells.Item(136, 10).Value = 109990.336
$ws.Cells.Item(136, 12).Value = 109990.336
$ws.Cells.Item(136, 14).Value = -120190.336
$ws.Cells.Item(137, 8).Value = 3810.1475
$ws.Cells.Item(137, 9).Value = 1814.3684
$ws.Cells.Item(137, 10).Value = 4713
$ws.Cells.Item(137, 11).Value = 5443.1052
$ws.Cells.Item(137, 12).Value = 14139
$ws.Cells.Item(137, 13).Value = -2893.1052
$ws.Cells.Item(137, 14).Value = -19239
$ws.Cells.Item(138, 8).Value = 5760.8438
$ws.Cells.Item(138, 9).Value = 3339.7778
$ws.Cells.Item(138, 10).Value = 6708.2173
$ws.Cells.Item(138, 11).Value = 10019.3334
$ws.Cells.Item(138, 12).Value = 20124.6519
$ws.Cells.Item(138, 13).Value = -4879.3334
$ws.Cells.Item(138, 14).Value = -30404.6519

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 2517.1428
$ws.Cells.Item(32, 9).Value = 2255.0286
$ws.Cells.Item(32, 10).Value = 3827.7144
$ws.Cells.Item(32, 11).Value = 2255.0286
$ws.Cells.Item(32, 12).Value = 3827.7144
$ws.Cells.Item(32, 13).Value = -1968.0286
$ws.Cells.Item(32, 14).Value = -4401.7144
$ws.Cells.Item(45, 8).Value = 50043184
$ws.Cells.Item(45, 9).Value = 53603.625
$ws.Cells.Item(45, 10).Value = 250001500
$ws.Cells.Item(45, 11).Value = 53603.625
$ws.Cells.Item(45, 12).Value = 250001500
$ws.Cells.Item(45, 13).Value = -53226.625
$ws.Cells.Item(45, 14).Value = -250002254
$ws.Cells.Item(61, 8).Value = 8073.8823
$ws.Cells.Item(61, 9).Value = 6648.6665
$ws.Cells.Item(61, 11).Value = 6648.6665
$ws.Cells.Item(61, 13).Value = -6436.6665
$ws.Cells.Item(74, 8).Value = 34483936
$ws.Cells.Item(74, 9).Value = 862.9524
$ws.Cells.Item(74, 10).Value = 125002000
$ws.Cells.Item(74, 11).Value = 862.9524
$ws.Cells.Item(74, 12).Value = 125002000
$ws.Cells.Item(74, 13).Value = 11.04759999999999
$ws.Cells.Item(74, 14).Value = -125003748
$ws.Cells.Item(77, 8).Value = 34483936
$ws.Cells.Item(77, 9).Value = 862.9524
$ws.Cells.Item(77, 10).Value = 125002000
$ws.Cells.Item(77, 11).Value = 4314.762
$ws.Cells.Item(77, 12).Value = 625010000
$ws.Cells.Item(77, 13).Value = 53.23800000000028
$ws.Cells.Item(77, 14).Value = -625018736
$ws.Cells.Item(110, 8).Value = 1003.6667
$ws.Cells.Item(110, 9).Value = 1003.6667
$ws.Cells.Item(110, 11).Value = 1003.6667
$ws.Cells.Item(110, 13).Value = 1041.3333
$ws.Cells.Item(122, 8).Value = 3050.0386
$ws.Cells.Item(122, 9).Value = 2394.5557
$ws.Cells.Item(122, 11).Value = 7183.6671
$ws.Cells.Item(122, 13).Value = -4733.6671
$ws.Cells.Item(132, 8).Value = 4709.5835
$ws.Cells.Item(132, 9).Value = 2525.238
$ws.Cells.Item(132, 10).Value = 20000
$ws.Cells.Item(132, 11).Value = 7575.714
$ws.Cells.Item(132, 12).Value = 60000
$ws.Cells.Item(132, 13).Value = -5045.714
$ws.Cells.Item(132, 14).Value = -65060
$ws.Cells.Item(135, 8).Value = 109923
$ws.Cells.Item(135, 10).Value = 109923
$ws.Cells.Item(135, 12).Value = 109923
$ws.Cells.Item(135, 14).Value = -120063
$ws.Cells.Item(136, 8).Value = 8073.8823
$ws.Cells.Item(136, 9).Value = 6648.6665
$ws.Cells.Item(136, 11).Value = 19945.9995
$ws.Cells.Item(136, 13).Value = -17395.9995

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 4587.077
$ws.Cells.Item(86, 9).Value = 4133.3
$ws.Cells.Item(86, 10).Value = 6099.6665
$ws.Cells.Item(86, 11).Value = 4133.3
$ws.Cells.Item(86, 12).Value = 6099.6665
$ws.Cells.Item(86, 13).Value = -3010.3
$ws.Cells.Item(86, 14).Value = -8345.666499999999
$ws.Cells.Item(89, 8).Value = 4587.077
$ws.Cells.Item(89, 9).Value = 4133.3
$ws.Cells.Item(89, 10).Value = 6099.6665
$ws.Cells.Item(89, 11).Value = 20666.5
$ws.Cells.Item(89, 12).Value = 30498.3325
$ws.Cells.Item(89, 13).Value = -15050.5
$ws.Cells.Item(89, 14).Value = -41730.3325
$ws.Cells.Item(107, 8).Value = 2653992.2
$ws.Cells.Item(107, 9).Value = 3078238
$ws.Cells.Item(107, 10).Value = 2455.75
$ws.Cells.Item(107, 11).Value = 3078238
$ws.Cells.Item(107, 12).Value = 2455.75
$ws.Cells.Item(107, 13).Value = -3076318
$ws.Cells.Item(107, 14).Value = -6295.75
$ws.Cells.Item(134, 8).Value = 2184.1667
$ws.Cells.Item(134, 10).Value = 2999.5
$ws.Cells.Item(134, 12).Value = 8998.5
$ws.Cells.Item(134, 14).Value = -14068.5
$ws.Cells.Item(135, 8).Value = 109997
$ws.Cells.Item(135, 10).Value = 109997
$ws.Cells.Item(135, 12).Value = 109997
$ws.Cells.Item(135, 14).Value = -120137
$ws.Cells.Item(140, 8).Value = 59996.5
$ws.Cells.Item(140, 10).Value = 0
$ws.Cells.Item(140, 12).Value = 0
$ws.Cells.Item(140, 14).ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 3077.554
$ws.Cells.Item(31, 9).Value = 2631.6
$ws.Cells.Item(31, 10).Value = 3147.2344
$ws.Cells.Item(31, 11).Value = 2631.6
$ws.Cells.Item(31, 12).Value = 3147.2344
$ws.Cells.Item(31, 13).Value = -2336.6
$ws.Cells.Item(31, 14).Value = -3737.2344
$ws.Cells.Item(34, 8).Value = 3077.554
$ws.Cells.Item(34, 9).Value = 2631.6
$ws.Cells.Item(34, 10).Value = 3147.2344
$ws.Cells.Item(34, 11).Value = 2631.6
$ws.Cells.Item(34, 12).Value = 3147.2344
$ws.Cells.Item(34, 13).Value = -2429.6
$ws.Cells.Item(34, 14).Value = -3551.2344
$ws.Cells.Item(58, 8).Value = 2471.2222
$ws.Cells.Item(58, 9).Value = 1167.1428
$ws.Cells.Item(58, 10).Value = 3301.0908
$ws.Cells.Item(58, 11).Value = 1167.1428
$ws.Cells.Item(58, 12).Value = 3301.0908
$ws.Cells.Item(58, 13).Value = -964.1428000000001
$ws.Cells.Item(58, 14).Value = -3707.0908
$ws.Cells.Item(105, 8).Value = 3002.5
$ws.Cells.Item(105, 9).Value = 2336.6667
$ws.Cells.Item(105, 11).Value = 2336.6667
$ws.Cells.Item(105, 13).Value = -589.6667000000002
$ws.Cells.Item(127, 8).Value = 33198.6
$ws.Cells.Item(127, 10).Value = 33198.6
$ws.Cells.Item(127, 12).Value = 33198.6
$ws.Cells.Item(127, 14).Value = -43118.6
$ws.Cells.Item(132, 8).Value = 3122.4285
$ws.Cells.Item(132, 9).Value = 2683.8696
$ws.Cells.Item(132, 10).Value = 5139.8
$ws.Cells.Item(132, 11).Value = 8051.6088
$ws.Cells.Item(132, 12).Value = 15419.4
$ws.Cells.Item(132, 13).Value = -5521.6088
$ws.Cells.Item(132, 14).Value = -20479.4
$ws.Cells.Item(133, 8).Value = 80872.375
$ws.Cells.Item(133, 10).Value = 81140
$ws.Cells.Item(133, 12).Value = 81140
$ws.Cells.Item(133, 14).Value = -86200
$ws.Cells.Item(134, 8).Value = 4035.5
$ws.Cells.Item(134, 9).Value = 3936.5557
$ws.Cells.Item(134, 10).Value = 4569.8
$ws.Cells.Item(134, 11).Value = 11809.6671
$ws.Cells.Item(134, 12).Value = 13709.4
$ws.Cells.Item(134, 13).Value = -9274.667099999999
$ws.Cells.Item(134, 14).Value = -18779.4
$ws.Cells.Item(135, 8).Value = 89997.664
$ws.Cells.Item(135, 10).Value = 89997.664
$ws.Cells.Item(135, 12).Value = 89997.664
$ws.Cells.Item(135, 14).Value = -100137.664
$ws.Cells.Item(136, 8).Value = 2471.2222
$ws.Cells.Item(136, 9).Value = 1167.1428
$ws.Cells.Item(136, 10).Value = 3301.0908
$ws.Cells.Item(136, 11).Value = 3501.4284
$ws.Cells.Item(136, 12).Value = 9903.2724
$ws.Cells.Item(136, 13).Value = -951.4284000000002
$ws.Cells.Item(136, 14).Value = -15003.2724
$ws.Cells.Item(137, 8).Value = 49999.5
$ws.Cells.Item(137, 10).Value = 49999.5
$ws.Cells.Item(137, 12).Value = 49999.5
$ws.Cells.Item(137, 14).Value = -60199.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(2, 8).Value = 287.5
$ws.Cells.Item(2, 9).Value = 210
$ws.Cells.Item(2, 11).Value = 1260
$ws.Cells.Item(2, 13).Value = -1147
$ws.Cells.Item(5, 8).Value = 1014.05554
$ws.Cells.Item(5, 10).Value = 1018.0714
$ws.Cells.Item(5, 12).Value = 3054.2142
$ws.Cells.Item(5, 14).Value = -3278.2142
$ws.Cells.Item(33, 8).Value = 177.8
$ws.Cells.Item(33, 9).Value = 274.5
$ws.Cells.Item(33, 10).Value = 113.333336
$ws.Cells.Item(33, 11).Value = 1647
$ws.Cells.Item(33, 12).Value = 680.000016
$ws.Cells.Item(33, 13).Value = -1364
$ws.Cells.Item(33, 14).Value = -1246.000016
$ws.Cells.Item(56, 8).Value = 8681.571
$ws.Cells.Item(56, 9).Value = 8681.571
$ws.Cells.Item(56, 11).Value = 8681.571
$ws.Cells.Item(56, 13).Value = -8151.571
$ws.Cells.Item(68, 8).Value = 1926651.4
$ws.Cells.Item(68, 10).Value = 2504012.5
$ws.Cells.Item(68, 12).Value = 7512037.5
$ws.Cells.Item(68, 14).Value = -7513659.5
$ws.Cells.Item(71, 8).Value = 1926651.4
$ws.Cells.Item(71, 10).Value = 2504012.5
$ws.Cells.Item(71, 12).Value = 22536112.5
$ws.Cells.Item(71, 14).Value = -22544224.5
$ws.Cells.Item(108, 8).Value = 149
$ws.Cells.Item(108, 9).Value = 149
$ws.Cells.Item(108, 11).Value = 447
$ws.Cells.Item(108, 13).Value = 2433
$ws.Cells.Item(129, 8).Value = 56820.2
$ws.Cells.Item(129, 10).Value = 63000.223
$ws.Cells.Item(129, 12).Value = 189000.669
$ws.Cells.Item(129, 14).Value = -199000.669
$ws.Cells.Item(132, 8).Value = 3252.5
$ws.Cells.Item(132, 9).Value = 0
$ws.Cells.Item(132, 10).Value = 3252.5
$ws.Cells.Item(132, 11).Value = 0
$ws.Cells.Item(132, 12).Value = 29272.5
$ws.Cells.Item(132, 13).ClearContents()
$ws.Cells.Item(132, 14).Value = -34332.5
$ws.Cells.Item(133, 8).Value = 0
$ws.Cells.Item(133, 9).Value = 0
$ws.Cells.Item(133, 11).Value = 0
$ws.Cells.Item(133, 13).ClearContents()
$ws.Cells.Item(135, 8).Value = 1014.05554
$ws.Cells.Item(135, 10).Value = 1018.0714
$ws.Cells.Item(135, 12).Value = 9162.642600000001
$ws.Cells.Item(135, 14).Value = -14232.6426
$ws.Cells.Item(139, 8).Value = 11750
$ws.Cells.Item(139, 10).Value = 6000
$ws.Cells.Item(139, 12).Value = 18000
$ws.Cells.Item(139, 14).Value = -28280

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 25004842
$ws.Cells.Item(70, 9).Value = 45458340
$ws.Cells.Item(70, 11).Value = 45458340
$ws.Cells.Item(70, 13).Value = -45458070
$ws.Cells.Item(73, 8).Value = 25004842
$ws.Cells.Item(73, 9).Value = 45458340
$ws.Cells.Item(73, 11).Value = 45458340
$ws.Cells.Item(73, 13).Value = -45457404
$ws.Cells.Item(97, 8).Value = 849.4815
$ws.Cells.Item(97, 9).Value = 726.2
$ws.Cells.Item(97, 11).Value = 726.2
$ws.Cells.Item(97, 13).Value = -230.2
$ws.Cells.Item(107, 8).Value = 8038.6
$ws.Cells.Item(107, 9).Value = 7398
$ws.Cells.Item(107, 11).Value = 7398
$ws.Cells.Item(107, 13).Value = -5478
$ws.Cells.Item(113, 8).Value = 250006560
$ws.Cells.Item(113, 9).Value = 187507500
$ws.Cells.Item(113, 11).Value = 187507500
$ws.Cells.Item(113, 13).Value = -187505330
$ws.Cells.Item(122, 8).Value = 100004130
$ws.Cells.Item(122, 9).Value = 4295.6665
$ws.Cells.Item(122, 10).Value = 250003870
$ws.Cells.Item(122, 11).Value = 12886.9995
$ws.Cells.Item(122, 12).Value = 750011610
$ws.Cells.Item(122, 13).Value = -10436.9995
$ws.Cells.Item(122, 14).Value = -750016510
$ws.Cells.Item(126, 8).Value = 55567212
$ws.Cells.Item(126, 9).Value = 83343150
$ws.Cells.Item(126, 10).Value = 15333
$ws.Cells.Item(126, 11).Value = 250029450
$ws.Cells.Item(126, 12).Value = 45999
$ws.Cells.Item(126, 13).Value = -250026980
$ws.Cells.Item(126, 14).Value = -50939
$ws.Cells.Item(132, 8).Value = 2518.1333
$ws.Cells.Item(132, 9).Value = 2219.111
$ws.Cells.Item(132, 10).Value = 2966.6667
$ws.Cells.Item(132, 11).Value = 6657.333
$ws.Cells.Item(132, 12).Value = 8900.000100000001
$ws.Cells.Item(132, 13).Value = -4127.333
$ws.Cells.Item(132, 14).Value = -13960.0001
$ws.Cells.Item(140, 8).Value = 74817.09
$ws.Cells.Item(140, 10).Value = 74817.09
$ws.Cells.Item(140, 12).Value = 74817.09
$ws.Cells.Item(140, 14).Value = -85177.09

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 10000
$ws.Cells.Item(7, 9).Value = 10000
$ws.Cells.Item(7, 10).Value = 0
$ws.Cells.Item(7, 11).Value = 10000
$ws.Cells.Item(7, 12).Value = 0
$ws.Cells.Item(7, 13).Value = -9888
$ws.Cells.Item(7, 14).ClearContents()
$ws.Cells.Item(16, 8).Value = 1457.6
$ws.Cells.Item(16, 9).Value = 1453.25
$ws.Cells.Item(16, 11).Value = 1453.25
$ws.Cells.Item(16, 13).Value = -1283.25
$ws.Cells.Item(40, 8).Value = 37667.086
$ws.Cells.Item(40, 9).Value = 45462.43
$ws.Cells.Item(40, 11).Value = 45462.43
$ws.Cells.Item(40, 13).Value = -45326.43
$ws.Cells.Item(55, 8).Value = 1434.45
$ws.Cells.Item(55, 9).Value = 1810.7693
$ws.Cells.Item(55, 11).Value = 1810.7693
$ws.Cells.Item(55, 13).Value = -1637.7693
$ws.Cells.Item(82, 8).Value = 3999.4
$ws.Cells.Item(82, 10).Value = 4998.5
$ws.Cells.Item(82, 12).Value = 4998.5
$ws.Cells.Item(82, 14).Value = -5720.5
$ws.Cells.Item(85, 8).Value = 3999.4
$ws.Cells.Item(85, 10).Value = 4998.5
$ws.Cells.Item(85, 12).Value = 4998.5
$ws.Cells.Item(85, 14).Value = -7494.5
$ws.Cells.Item(100, 8).Value = 3126.4285
$ws.Cells.Item(100, 9).Value = 2746.5
$ws.Cells.Item(100, 10).Value = 3633
$ws.Cells.Item(100, 11).Value = 2746.5
$ws.Cells.Item(100, 12).Value = 3633
$ws.Cells.Item(100, 13).Value = -2205.5
$ws.Cells.Item(100, 14).Value = -4715
$ws.Cells.Item(122, 8).Value = 7499.1665
$ws.Cells.Item(122, 9).Value = 6999
$ws.Cells.Item(122, 11).Value = 20997
$ws.Cells.Item(122, 13).Value = -18547
$ws.Cells.Item(126, 8).Value = 10000
$ws.Cells.Item(126, 9).Value = 10000
$ws.Cells.Item(126, 10).Value = 0
$ws.Cells.Item(126, 11).Value = 30000
$ws.Cells.Item(126, 12).Value = 0
$ws.Cells.Item(126, 13).Value = -27530
$ws.Cells.Item(126, 14).ClearContents()
$ws.Cells.Item(140, 8).Value = 104185.44
$ws.Cells.Item(140, 10).Value = 104185.44
$ws.Cells.Item(140, 12).Value = 104185.44
$ws.Cells.Item(140, 14).Value = -114545.44

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(69, 8).Value = 0
$ws.Cells.Item(69, 10).Value = 0
$ws.Cells.Item(69, 12).Value = 0
$ws.Cells.Item(69, 14).ClearContents()
$ws.Cells.Item(72, 8).Value = 0
$ws.Cells.Item(72, 10).Value = 0
$ws.Cells.Item(72, 12).Value = 0
$ws.Cells.Item(72, 14).ClearContents()
$ws.Cells.Item(96, 8).Value = 1098.1177
$ws.Cells.Item(96, 9).Value = 1114
$ws.Cells.Item(96, 10).Value = 1089.4546
$ws.Cells.Item(96, 11).Value = 1114
$ws.Cells.Item(96, 12).Value = 1089.4546
$ws.Cells.Item(96, 13).Value = 259
$ws.Cells.Item(96, 14).Value = -3835.4546
$ws.Cells.Item(122, 8).Value = 31253068
$ws.Cells.Item(122, 9).Value = 2134.75
$ws.Cells.Item(122, 10).Value = 62504000
$ws.Cells.Item(122, 11).Value = 6404.25
$ws.Cells.Item(122, 12).Value = 187512000
$ws.Cells.Item(122, 13).Value = -3954.25
$ws.Cells.Item(122, 14).Value = -187516900
$ws.Cells.Item(126, 8).Value = 6185.6816
$ws.Cells.Item(126, 10).Value = 3520
$ws.Cells.Item(126, 12).Value = 10560
$ws.Cells.Item(126, 14).Value = -15500
$ws.Cells.Item(132, 8).Value = 7041.9375
$ws.Cells.Item(132, 9).Value = 6913.96
$ws.Cells.Item(132, 10).Value = 7499
$ws.Cells.Item(132, 11).Value = 20741.88
$ws.Cells.Item(132, 12).Value = 22497
$ws.Cells.Item(132, 13).Value = -18211.88
$ws.Cells.Item(132, 14).Value = -27557
$ws.Cells.Item(136, 8).Value = 18534556
$ws.Cells.Item(136, 9).Value = 20849000
$ws.Cells.Item(136, 11).Value = 62547000
$ws.Cells.Item(136, 13).Value = -62544450
$ws.Cells.Item(138, 8).Value = 191666.67
$ws.Cells.Item(138, 10).Value = 191666.67
$ws.Cells.Item(138, 12).Value = 191666.67
$ws.Cells.Item(138, 14).Value = -201946.67
